$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add summary rows (order chosen to match original authoring / shared-string order)
$ws.Range("F14").Value = "Average:"
$ws.Range("G14").Formula = "=AVERAGE(G2:G12)"

$ws.Range("F15").Value = "Std Deviation:"
$ws.Range("G15").Formula = "=STDEV.P(G2:G12)"

$ws.Range("F17").Value = "Low Average"
$ws.Range("G17").Formula = "=AVERAGE(G2,G4,G7,G8,G10,G11)"

$ws.Range("F18").Value = "Low Std. Dev"
$ws.Range("G18").Formula = "=STDEV.S(G2,G4,G7,G10,G11)"

$ws.Range("F21").Value = "High Std. Dev"
$ws.Range("G21").Formula = "=STDEV.S(G5,G3,G6,G9,G12)"

$ws.Range("F20").Value = "High Average"
$ws.Range("G20").Formula = "=AVERAGE(G3,G5,G6,G9,G12)"

# 2. Update the "hops" header to "hops to target"
$ws.Range("B1").Value = "hops to target"

# 3. Add new "time/hop (corrected)" column G
$ws.Range("G1").Value = "time/hop (corrected)"
$ws.Range("G2").Formula = "=C2/(B2*2)"
$ws.Range("G3:G12").Formula = "=C3/(B3*2)"

# 4. Autofit the new/changed columns
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(7).AutoFit() | Out-Null

# 5. Reposition the chart: it shifts from anchor H2:P21 to I2:R21 (one column right,
# two columns wider) now that the new columns push it over.
$co = $ws.ChartObjects().Item(1)

$emuPerPt = 12700.0

$left = 0.0
for ($c = 1; $c -le 8; $c++) { $left += $ws.Columns.Item($c).Width }
$left += 1 / $emuPerPt

$top = 14.4 * 1
$top += 1 / $emuPerPt

$right = 0.0
for ($c = 1; $c -le 17; $c++) { $right += $ws.Columns.Item($c).Width }
$right += 1 / $emuPerPt

$bottom = 14.4 * 20
$bottom += 1 / $emuPerPt

$co.Left = $left
$co.Top = $top
$co.Width = $right - $left
$co.Height = $bottom - $top

# 6. Update selection
$ws.Range("R12").Select()
